$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 8
$ws.Range("H8").Value = 46.57143
$ws.Range("I8").Value = 46.57143
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 139.71429
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -0.7142900000000054
$ws.Range("N8").ClearContents()
# row 15
$ws.Range("H15").Value = 4127.1064
$ws.Range("I15").Value = 4127.1064
$ws.Range("K15").Value = 12381.3192
$ws.Range("M15").Value = -12212.3192
# row 18
$ws.Range("H18").Value = 1275.8
$ws.Range("I18").Value = 1275.8
$ws.Range("K18").Value = 1275.8
$ws.Range("M18").Value = -991.8
# row 64
$ws.Range("H64").Value = 4771.4287
$ws.Range("I64").Value = 3800
$ws.Range("J64").Value = 4933.3335
$ws.Range("K64").Value = 3800
$ws.Range("L64").Value = 4933.3335
$ws.Range("M64").Value = -3552
$ws.Range("N64").Value = -5429.3335
# row 67
$ws.Range("H67").Value = 4771.4287
$ws.Range("I67").Value = 3800
$ws.Range("J67").Value = 4933.3335
$ws.Range("K67").Value = 3800
$ws.Range("L67").Value = 4933.3335
$ws.Range("M67").Value = -2942
$ws.Range("N67").Value = -6649.3335
# row 69
$ws.Range("H69").Value = 3015
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3015
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 9045
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -10793
# row 72
$ws.Range("H72").Value = 3015
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3015
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 27135
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -35871
# row 76
$ws.Range("H76").Value = 3465.6667
$ws.Range("I76").Value = 3400
$ws.Range("K76").Value = 3400
$ws.Range("M76").Value = -3085
# row 79
$ws.Range("H79").Value = 3465.6667
$ws.Range("I79").Value = 3400
$ws.Range("K79").Value = 3400
$ws.Range("M79").Value = -2308
# row 132
$ws.Range("H132").Value = 732861.6
$ws.Range("I132").Value = 1359.2142
$ws.Range("J132").Value = 4456874
$ws.Range("K132").Value = 4077.6426
$ws.Range("L132").Value = 13370622
$ws.Range("M132").Value = -1547.6426
$ws.Range("N132").Value = -13375682

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 2770.0625
$ws.Range("I2").Value = 2827.75
$ws.Range("K2").Value = 2827.75
$ws.Range("M2").Value = -2714.75
# row 24
$ws.Range("H24").Value = 25870.572
$ws.Range("J24").Value = 25870.572
$ws.Range("L24").Value = 25870.572
$ws.Range("N24").Value = -26618.572
# row 32
$ws.Range("H32").Value = 24206.586
$ws.Range("I32").Value = 39162.547
$ws.Range("J32").Value = 7034.926
$ws.Range("K32").Value = 39162.547
$ws.Range("L32").Value = 7034.926
$ws.Range("M32").Value = -38875.547
$ws.Range("N32").Value = -7608.926
# row 61
$ws.Range("H61").Value = 52737824
$ws.Range("I61").Value = 77000660
$ws.Range("J61").Value = 168366.67
$ws.Range("K61").Value = 77000660
$ws.Range("L61").Value = 168366.67
$ws.Range("M61").Value = -77000448
$ws.Range("N61").Value = -168790.67
# row 100
$ws.Range("H100").Value = 25870.572
$ws.Range("J100").Value = 25870.572
$ws.Range("L100").Value = 25870.572
$ws.Range("N100").Value = -28034.572
# row 110
$ws.Range("H110").Value = 770675.0600000001
$ws.Range("I110").Value = 1429711.4
$ws.Range("J110").Value = 1799.3334
$ws.Range("K110").Value = 1429711.4
$ws.Range("L110").Value = 1799.3334
$ws.Range("M110").Value = -1427666.4
$ws.Range("N110").Value = -5889.3334
# row 116
$ws.Range("H116").Value = 2770.0625
$ws.Range("I116").Value = 2827.75
$ws.Range("K116").Value = 2827.75
$ws.Range("M116").Value = -533.75
# row 136
$ws.Range("H136").Value = 52737824
$ws.Range("I136").Value = 77000660
$ws.Range("J136").Value = 168366.67
$ws.Range("K136").Value = 231001980
$ws.Range("L136").Value = 505100.01
$ws.Range("M136").Value = -230999430
$ws.Range("N136").Value = -510200.01

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 2770.0625
$ws.Range("I3").Value = 2827.75
$ws.Range("K3").Value = 2827.75
$ws.Range("M3").Value = -2713.75
# row 20
$ws.Range("H20").Value = 1710.7
$ws.Range("I20").Value = 1514
$ws.Range("J20").Value = 2169.6667
$ws.Range("K20").Value = 1514
$ws.Range("L20").Value = 2169.6667
$ws.Range("M20").Value = -1267
$ws.Range("N20").Value = -2663.6667
# row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# row 105
$ws.Range("H105").Value = 33335634
$ws.Range("I105").Value = 45456876
$ws.Range("J105").Value = 2225
$ws.Range("K105").Value = 45456876
$ws.Range("L105").Value = 2225
$ws.Range("M105").Value = -45455129
$ws.Range("N105").Value = -5719
# row 132
$ws.Range("H132").Value = 48890
$ws.Range("J132").Value = 48890
$ws.Range("L132").Value = 48890
$ws.Range("N132").Value = -59010
# row 134
$ws.Range("H134").Value = 2219.1843
$ws.Range("I134").Value = 1222.1818
$ws.Range("J134").Value = 3590.0625
$ws.Range("K134").Value = 3666.5454
$ws.Range("L134").Value = 10770.1875
$ws.Range("M134").Value = -1131.5454
$ws.Range("N134").Value = -15840.1875

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 62
$ws.Range("H62").Value = 2501.2
$ws.Range("I62").Value = 2125
$ws.Range("J62").Value = 4006
$ws.Range("K62").Value = 2125
$ws.Range("L62").Value = 4006
$ws.Range("M62").Value = -1501
$ws.Range("N62").Value = -5254
# row 65
$ws.Range("H65").Value = 2501.2
$ws.Range("I65").Value = 2125
$ws.Range("J65").Value = 4006
$ws.Range("K65").Value = 10625
$ws.Range("L65").Value = 20030
$ws.Range("M65").Value = -7505
$ws.Range("N65").Value = -26270
# row 134
$ws.Range("H134").Value = 20654.232
$ws.Range("I134").Value = 1418.8975
$ws.Range("J134").Value = 64782.35
$ws.Range("K134").Value = 4256.6925
$ws.Range("L134").Value = 194347.05
$ws.Range("M134").Value = -1721.6925
$ws.Range("N134").Value = -199417.05

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 56
$ws.Range("H56").Value = 154861.05
$ws.Range("I56").Value = 154861.05
$ws.Range("K56").Value = 154861.05
$ws.Range("M56").Value = -154331.05

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 10
$ws.Range("H10").Value = 500
$ws.Range("I10").Value = 500
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 500
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -360
$ws.Range("N10").ClearContents()
# row 69
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
# row 72
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
# row 136
$ws.Range("H136").Value = 145298.94
$ws.Range("I136").Value = 159497.72
$ws.Range("J136").Value = 132875
$ws.Range("K136").Value = 478493.16
$ws.Range("L136").Value = 398625
$ws.Range("M136").Value = -475943.16
$ws.Range("N136").Value = -403725

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
# row 14
$ws.Range("H14").Value = 13999.833
$ws.Range("J14").Value = 10800
$ws.Range("L14").Value = 10800
$ws.Range("N14").Value = -11136
# row 132
$ws.Range("H132").Value = 56151.07
$ws.Range("I132").Value = 40786
$ws.Range("J132").Value = 184193.33
$ws.Range("K132").Value = 122358
$ws.Range("L132").Value = 552579.99
$ws.Range("M132").Value = -119828
$ws.Range("N132").Value = -557639.99
# row 136
$ws.Range("H136").Value = 47048.977
$ws.Range("I136").Value = 34258.5
$ws.Range("J136").Value = 74457.14
$ws.Range("K136").Value = 102775.5
$ws.Range("L136").Value = 223371.42
$ws.Range("M136").Value = -100225.5
$ws.Range("N136").Value = -228471.42
